# Apply the "LOG" cleanup + test-case text fix described in the commit:
#   - Remove the per-row "LOG" log-message entries in column E (rows 2-5),
#     leaving the cells blank but keeping their existing style/border.
#   - Drop the trailing period from the Twitter row's expected text,
#     "Join Twitter today." -> "Join Twitter today".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# E2:E5 held the literal "LOG" placeholder; clear the values but keep formatting.
$ws.Range("E2:E5").ClearContents()

# B4 ("Navigation_Url" row for twitter.com) had a trailing period; remove it.
$ws.Range("B4").Value = "Join Twitter today"
